$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (Price + Volume(1h)) and three row swaps
# (Hedera <-> PancakeSwap, MXToken <-> ARBITRUM, TrustWalletToken <-> Quant)
# Values in column D that look numeric are prefixed with a literal apostrophe
# so Excel stores them as text (matching the source data format), not numbers.
$updates = @{
    'D2' = '25.840.99'
    'E2' = '  -1.23%  '
    'D3' = '1.635.95'
    'E3' = '  -1.38%  '
    'D4' = '''1.003'
    'E4' = '  -0.27%  '
    'D5' = '''215.49'
    'E5' = '  -1.38%  '
    'D6' = '''0.5020'
    'E6' = '  -2.64%  '
    'D7' = '''1.004'
    'E7' = '  -0.24%  '
    'D8' = '''0.2568'
    'E8' = '  -0.82%  '
    'D9' = '''0.06421'
    'E9' = '  -0.45%  '
    'D10' = '''19.60'
    'E10' = '  -1.60%  '
    'D11' = '''0.07667'
    'E11' = '  -1.60%  '
    'D12' = '1.636.10'
    'E12' = '  -1.55%  '
    'D13' = '''4.241'
    'E13' = '  -1.29%  '
    'D14' = '1.861.63'
    'E14' = '  -1.35%  '
    'D15' = '''0.5456'
    'E15' = '  -1.93%  '
    'D16' = '0.0₅7931'
    'E16' = '  -1.63%  '
    'D17' = '''63.48'
    'E17' = '  -1.21%  '
    'D18' = '25.864.67'
    'E18' = '  -1.25%  '
    'D19' = '''1.003'
    'E19' = '  -0.23%  '
    'D20' = '''203.13'
    'E20' = '  -4.25%  '
    'D21' = '''4.308'
    'E21' = '  -2.37%  '
    'D22' = '''9.952'
    'E22' = '  -0.88%  '
    'D23' = '''5.997'
    'E23' = '  +0.38%  '
    'D24' = '''1.004'
    'E24' = '  -0.22%  '
    'D25' = '''1.935'
    'E25' = '  +10.03%  '
    'D26' = '''141.56'
    'E26' = '  -2.07%  '
    'D27' = '''0.1147'
    'E27' = '  -1.47%  '
    'D28' = '''15.69'
    'E28' = '  -0.79%  '
    'D29' = '''6.707'
    'E29' = '  -3.87%  '
    'B30' = 'PancakeSwap'
    'C30' = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
    'D30' = '''1.240'
    'E30' = '  -1.19%  '
    'B31' = 'Hedera'
    'C31' = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
    'D31' = '''0.05002'
    'E31' = '  -5.25%  '
    'D32' = '''3.261'
    'E32' = '  -3.13%  '
    'D33' = '''3.188'
    'E33' = '  -0.95%  '
    'D34' = '''1.531'
    'E34' = '  -2.57%  '
    'D35' = '''2.351'
    'E35' = '  -0.90%  '
    'D36' = '1.172.49'
    'E36' = '  +0.47%  '
    'B37' = 'ARBITRUM'
    'C37' = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
    'D37' = '''0.8931'
    'E37' = '  -3.87%  '
    'B38' = 'MXToken'
    'C38' = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
    'D38' = '''2.623'
    'E38' = '  -5.08%  '
    'D39' = '''0.5561'
    'E39' = '  -1.73%  '
    'D40' = '''0.01560'
    'E40' = '  -2.11%  '
    'D41' = '''2.556'
    'E41' = '  -0.24%  '
    'D42' = '''1.003'
    'E42' = '  -0.23%  '
    'D43' = '''5.645'
    'E43' = '  -0.81%  '
    'B44' = 'Quant'
    'C44' = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
    'D44' = '''99.69'
    'E44' = '  -0.63%  '
    'B45' = 'TrustWalletToken'
    'C45' = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
    'D45' = '''0.8042'
    'E45' = '  -4.77%  '
    'D46' = '1.775.40'
    'E46' = '  -1.21%  '
    'D47' = '0.0₈108'
    'E47' = '  -4.14%  '
    'D48' = '''0.4511'
    'E48' = '  -0.56%  '
    'D49' = '''1.006'
    'E49' = '  +0.12%  '
    'D50' = '''54.89'
    'E50' = '  -1.75%  '
    'D51' = '''0.05036'
    'E51' = '  -0.39%  '
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

Write-Host "Applied" $updates.Count "cell updates"
